$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.2994946666666667
$ws.Range("H2").Value = 0.8984840000000001
$ws.Range("I2").Value = 0.4989451716962827
$ws.Range("J2").Value = 0.4989451716962828
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.1145763333333333
$ws.Range("N2").Value = 0.343729
$ws.Range("O2").Value = 0.006557053879060051
$ws.Range("P2").Value = 0.006557053879060051
$ws.Range("Q2").Value = 0.03431500075955556
$ws.Range("R2").Value = 0.308835006836
$ws.Range("S2").Value = 0.003271610373509393
$ws.Range("T2").Value = 0.003271610373509394
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.2994946666666667
$ws.Range("H3").Value = 0.8984840000000001
$ws.Range("I3").Value = 0.4989451716962827
$ws.Range("J3").Value = 0.4989451716962828
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.390663666666667
$ws.Range("N3").Value = 28.171991
$ws.Range("O3").Value = 0.5374154140831726
$ws.Range("P3").Value = 0.5374154140831726
$ws.Range("Q3").Value = 2.812453684627112
$ws.Range("R3").Value = 25.312083161644
$ws.Range("S3").Value = 0.2681408260519574
$ws.Range("T3").Value = 0.2681408260519574
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.2994946666666667
$ws.Range("H4").Value = 0.8984840000000001
$ws.Range("I4").Value = 0.4989451716962827
$ws.Range("J4").Value = 0.4989451716962828
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.968512
$ws.Range("N4").Value = 23.905536
$ws.Range("O4").Value = 0.4560275320377672
$ws.Range("P4").Value = 0.4560275320377672
$ws.Range("Q4").Value = 2.386526845269334
$ws.Range("R4").Value = 21.478741607424
$ws.Range("S4").Value = 0.2275327352708159
$ws.Range("T4").Value = 0.2275327352708159
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.300761
$ws.Range("H5").Value = 0.9022829999999999
$ws.Range("I5").Value = 0.5010548283037172
$ws.Range("J5").Value = 0.5010548283037172
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.1145763333333333
$ws.Range("N5").Value = 0.343729
$ws.Range("O5").Value = 0.006557053879060051
$ws.Range("P5").Value = 0.006557053879060051
$ws.Range("Q5").Value = 0.03446009258966666
$ws.Range("R5").Value = 0.3101408333069999
$ws.Range("S5").Value = 0.003285443505550656
$ws.Range("T5").Value = 0.003285443505550656
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.300761
$ws.Range("H6").Value = 0.9022829999999999
$ws.Range("I6").Value = 0.5010548283037172
$ws.Range("J6").Value = 0.5010548283037172
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.390663666666667
$ws.Range("N6").Value = 28.171991
$ws.Range("O6").Value = 0.5374154140831726
$ws.Range("P6").Value = 0.5374154140831726
$ws.Range("Q6").Value = 2.824345395050333
$ws.Range("R6").Value = 25.419108555453
$ws.Range("S6").Value = 0.2692745880312151
$ws.Range("T6").Value = 0.2692745880312151
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.300761
$ws.Range("H7").Value = 0.9022829999999999
$ws.Range("I7").Value = 0.5010548283037172
$ws.Range("J7").Value = 0.5010548283037172
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.968512
$ws.Range("N7").Value = 23.905536
$ws.Range("O7").Value = 0.4560275320377672
$ws.Range("P7").Value = 0.4560275320377672
$ws.Range("Q7").Value = 2.396617637632
$ws.Range("R7").Value = 21.569558738688
$ws.Range("S7").Value = 0.2284947967669513
$ws.Range("T7").Value = 0.2284947967669513
